# Generate Report for Archive
#
# The "5641e55e-733f-4435-8657-1fdd4b40b28b" entry moves from the last row
# (row 7) of each status table up to the second data row (row 4), pushing
# the "875bc3c4...", "18952b98..." and "51f9abbe..." rows down by one.
# This touches the three worksheets (Overview, zh-cn, de-de) plus every
# hyperlink anchored on column B (Overview) / column A (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A..G, data rows 2..7
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = "5641e55e-733f-4435-8657-1fdd4b40b28b.md"
$ws1.Range("B4").Value = "e2e\5641e55e-733f-4435-8657-1fdd4b40b28b.md"
$ws1.Range("C4").Value = ".md"
$ws1.Range("D4").Value = ""
$ws1.Range("E4").Value = "Ready for handoff"
$ws1.Range("F4").Value = "Ready for handoff"
$ws1.Range("G4").Value = "2016-08-22 12:43:31"

$ws1.Range("A5").Value = "875bc3c4-bb3b-495f-85aa-7f2348317857.md"
$ws1.Range("B5").Value = "e2e\875bc3c4-bb3b-495f-85aa-7f2348317857.md"
$ws1.Range("C5").Value = ".md"
$ws1.Range("D5").Value = ""
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("G5").Value = "2016-08-22 12:42:54"

$ws1.Range("A6").Value = "18952b98-f01f-43a0-94af-4b31ebabd2c1.md"
$ws1.Range("B6").Value = "e2e\18952b98-f01f-43a0-94af-4b31ebabd2c1.md"
$ws1.Range("C6").Value = ".md"
$ws1.Range("D6").Value = ""
$ws1.Range("E6").Value = "Ready for handoff"
$ws1.Range("F6").Value = "Ready for handoff"
$ws1.Range("G6").Value = "2016-08-22 12:41:41"

$ws1.Range("A7").Value = "51f9abbe-2412-45cf-881e-4eecab8e723b.md"
$ws1.Range("B7").Value = "e2e\51f9abbe-2412-45cf-881e-4eecab8e723b.md"
$ws1.Range("C7").Value = ".md"
$ws1.Range("D7").Value = ""
$ws1.Range("E7").Value = "Ready for handoff"
$ws1.Range("F7").Value = "Ready for handoff"
$ws1.Range("G7").Value = "2016-08-22 12:43:46"

# Rebuild the hyperlinks on column B in the new row order. Deleting via a
# Range clears every hyperlink on the sheet, so it only needs doing once.
$ws1.Range("B2").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e889ed889379c5fb1135bd325feab0c58576d417/e2e/10fd8f40-bc55-457f-a985-0b9421aa9718.md", "", "", "e2e\10fd8f40-bc55-457f-a985-0b9421aa9718.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feeea8cd3eca2709645a869a651e11440e212c10/e2e/4bc103c6-5349-45d4-b84b-a40326524a31.md", "", "", "e2e\4bc103c6-5349-45d4-b84b-a40326524a31.md")
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/641656f3242baf9bb31575788873530ecd89051c/e2e/5641e55e-733f-4435-8657-1fdd4b40b28b.md", "", "", "e2e\5641e55e-733f-4435-8657-1fdd4b40b28b.md")
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feeea8cd3eca2709645a869a651e11440e212c10/e2e/875bc3c4-bb3b-495f-85aa-7f2348317857.md", "", "", "e2e\875bc3c4-bb3b-495f-85aa-7f2348317857.md")
$ws1.Hyperlinks.Add($ws1.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a05054c3d64e8d8a6dc37908b2a02b70eaa4a8c/e2e/18952b98-f01f-43a0-94af-4b31ebabd2c1.md", "", "", "e2e\18952b98-f01f-43a0-94af-4b31ebabd2c1.md")
$ws1.Hyperlinks.Add($ws1.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c89a498d93849f38476278d46e1857ea680e5e98/e2e/51f9abbe-2412-45cf-881e-4eecab8e723b.md", "", "", "e2e\51f9abbe-2412-45cf-881e-4eecab8e723b.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A..P, data rows 2..7
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value = "5641e55e-733f-4435-8657-1fdd4b40b28b.md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("G4").Value = "5641e55e-733f-4435-8657-1fdd4b40b28b.9df24d998f65fc2f6d717eaa2563ba539f8bc04e.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-22 12:43:26"
$ws2.Range("K4").Value = "0001-01-01 00:00:00"

$ws2.Range("A5").Value = "875bc3c4-bb3b-495f-85aa-7f2348317857.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("G5").Value = "875bc3c4-bb3b-495f-85aa-7f2348317857.002bca59b1e4767c6feb098c2c21a1f3d2249538.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-22 12:42:49"
$ws2.Range("K5").Value = "0001-01-01 00:00:00"

$ws2.Range("A6").Value = "18952b98-f01f-43a0-94af-4b31ebabd2c1.md"
$ws2.Range("C6").Value = "Ready for handoff"
$ws2.Range("G6").Value = "18952b98-f01f-43a0-94af-4b31ebabd2c1.ca8633e8d73304cfb223849ae911854a18667172.zh-cn.xlf"
$ws2.Range("H6").Value = "2016-08-22 12:41:37"
$ws2.Range("K6").Value = "0001-01-01 00:00:00"

$ws2.Range("A7").Value = "51f9abbe-2412-45cf-881e-4eecab8e723b.md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("G7").Value = "51f9abbe-2412-45cf-881e-4eecab8e723b.3f45bc725158022bf7e73ef6b835ad259349faef.zh-cn.xlf"
$ws2.Range("H7").Value = "2016-08-22 12:43:42"
$ws2.Range("K7").Value = "0001-01-01 00:00:00"

$ws2.Range("A2").Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e889ed889379c5fb1135bd325feab0c58576d417/e2e/10fd8f40-bc55-457f-a985-0b9421aa9718.md", "", "", "10fd8f40-bc55-457f-a985-0b9421aa9718.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2ef431c76bdca0e840bc89dc44c03ec268ab8158/e2e/10fd8f40-bc55-457f-a985-0b9421aa9718.md", "", "", "10fd8f40-bc55-457f-a985-0b9421aa9718.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feeea8cd3eca2709645a869a651e11440e212c10/e2e/4bc103c6-5349-45d4-b84b-a40326524a31.md", "", "", "4bc103c6-5349-45d4-b84b-a40326524a31.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/641656f3242baf9bb31575788873530ecd89051c/e2e/5641e55e-733f-4435-8657-1fdd4b40b28b.md", "", "", "5641e55e-733f-4435-8657-1fdd4b40b28b.md")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feeea8cd3eca2709645a869a651e11440e212c10/e2e/875bc3c4-bb3b-495f-85aa-7f2348317857.md", "", "", "875bc3c4-bb3b-495f-85aa-7f2348317857.md")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a05054c3d64e8d8a6dc37908b2a02b70eaa4a8c/e2e/18952b98-f01f-43a0-94af-4b31ebabd2c1.md", "", "", "18952b98-f01f-43a0-94af-4b31ebabd2c1.md")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c89a498d93849f38476278d46e1857ea680e5e98/e2e/51f9abbe-2412-45cf-881e-4eecab8e723b.md", "", "", "51f9abbe-2412-45cf-881e-4eecab8e723b.md")

# ---------------------------------------------------------------------
# Sheet "de-de": columns A..P, data rows 2..7
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value = "5641e55e-733f-4435-8657-1fdd4b40b28b.md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("G4").Value = "5641e55e-733f-4435-8657-1fdd4b40b28b.9df24d998f65fc2f6d717eaa2563ba539f8bc04e.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-22 12:43:31"
$ws3.Range("K4").Value = "0001-01-01 00:00:00"

$ws3.Range("A5").Value = "875bc3c4-bb3b-495f-85aa-7f2348317857.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("G5").Value = "875bc3c4-bb3b-495f-85aa-7f2348317857.002bca59b1e4767c6feb098c2c21a1f3d2249538.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-22 12:42:54"
$ws3.Range("K5").Value = "0001-01-01 00:00:00"

$ws3.Range("A6").Value = "18952b98-f01f-43a0-94af-4b31ebabd2c1.md"
$ws3.Range("C6").Value = "Ready for handoff"
$ws3.Range("G6").Value = "18952b98-f01f-43a0-94af-4b31ebabd2c1.ca8633e8d73304cfb223849ae911854a18667172.de-de.xlf"
$ws3.Range("H6").Value = "2016-08-22 12:41:41"
$ws3.Range("K6").Value = "0001-01-01 00:00:00"

$ws3.Range("A7").Value = "51f9abbe-2412-45cf-881e-4eecab8e723b.md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("G7").Value = "51f9abbe-2412-45cf-881e-4eecab8e723b.3f45bc725158022bf7e73ef6b835ad259349faef.de-de.xlf"
$ws3.Range("H7").Value = "2016-08-22 12:43:46"
$ws3.Range("K7").Value = "0001-01-01 00:00:00"

$ws3.Range("A2").Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e889ed889379c5fb1135bd325feab0c58576d417/e2e/10fd8f40-bc55-457f-a985-0b9421aa9718.md", "", "", "10fd8f40-bc55-457f-a985-0b9421aa9718.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/59bbc1fab957d301d4c14b8b6bbd59aa4a1cbf9e/e2e/10fd8f40-bc55-457f-a985-0b9421aa9718.md", "", "", "10fd8f40-bc55-457f-a985-0b9421aa9718.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feeea8cd3eca2709645a869a651e11440e212c10/e2e/4bc103c6-5349-45d4-b84b-a40326524a31.md", "", "", "4bc103c6-5349-45d4-b84b-a40326524a31.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/641656f3242baf9bb31575788873530ecd89051c/e2e/5641e55e-733f-4435-8657-1fdd4b40b28b.md", "", "", "5641e55e-733f-4435-8657-1fdd4b40b28b.md")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/feeea8cd3eca2709645a869a651e11440e212c10/e2e/875bc3c4-bb3b-495f-85aa-7f2348317857.md", "", "", "875bc3c4-bb3b-495f-85aa-7f2348317857.md")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a05054c3d64e8d8a6dc37908b2a02b70eaa4a8c/e2e/18952b98-f01f-43a0-94af-4b31ebabd2c1.md", "", "", "18952b98-f01f-43a0-94af-4b31ebabd2c1.md")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c89a498d93849f38476278d46e1857ea680e5e98/e2e/51f9abbe-2412-45cf-881e-4eecab8e723b.md", "", "", "51f9abbe-2412-45cf-881e-4eecab8e723b.md")
